# Fill in the missing "Absent" values (column H) to complete the
# consolidated report. Absent = 1 - Real (column E) for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value  = 1
$ws.Range("H5").Value  = 0
$ws.Range("H9").Value  = 1
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
